$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new literal text value.
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (prices, percentages, hour counters) are
# stored as literal text instead of being auto-converted to numbers,
# matching the source data which keeps these values as text.
$cellUpdates = @(
    @("D2", "331.27"),
    @("E2", "0.80%"),
    @("G2", "16"),
    @("E3", "2.66%"),
    @("G3", "16"),
    @("D4", "5.738"),
    @("E4", "0.11%"),
    @("G4", "16"),
    @("D5", "0.08107"),
    @("E5", "-0.07%"),
    @("G5", "16"),
    @("D6", "8.674"),
    @("E6", "-0.24%"),
    @("G6", "16"),
    @("D7", "4.498"),
    @("E7", "-1.72%"),
    @("G7", "16"),
    @("D8", "1.985"),
    @("E8", "1.53%"),
    @("G8", "16"),
    @("D9", "2.928"),
    @("E9", "-2.38%"),
    @("G9", "16"),
    @("D10", "0.9265"),
    @("E10", "-1.87%"),
    @("G10", "16"),
    @("D11", "0.1282"),
    @("E11", "-0.70%"),
    @("G11", "16"),
    @("D12", "0.1956"),
    @("E12", "-2.25%"),
    @("G12", "16"),
    @("D13", "8.822"),
    @("E13", "15.87%"),
    @("G13", "16"),
    @("D14", "0.09212"),
    @("E14", "0.55%"),
    @("G14", "16"),
    @("D15", "0.03734"),
    @("E15", "7.72%"),
    @("G15", "16"),
    @("D16", "0.1049"),
    @("E16", "9.21%"),
    @("G16", "16"),
    @("D17", "0.001296"),
    @("E17", "-3.12%"),
    @("G17", "16"),
    @("D18", "0.006268"),
    @("E18", "2.64%"),
    @("G18", "16"),
    @("D19", "3.369"),
    @("G19", "16"),
    @("E20", "-1.01%"),
    @("G20", "16"),
    @("D21", "0.1378"),
    @("G21", "16"),
    @("D22", "0.2606"),
    @("E22", "7.61%"),
    @("G22", "16"),
    @("D23", "0.04421"),
    @("E23", "-0.52%"),
    @("G23", "16"),
    @("D24", "0.001253"),
    @("E24", "0.11%"),
    @("G24", "16"),
    @("D25", "0.004405"),
    @("E25", "0.82%"),
    @("G25", "16"),
    @("E26", "4.19%"),
    @("G26", "16"),
    @("G27", "16"),
    @("G28", "16"),
    @("G29", "16"),
    @("G30", "16"),
    @("G31", "16"),
    @("G32", "16"),
    @("G33", "16"),
    @("G34", "16"),
    @("G35", "16"),
    @("G36", "16"),
    @("G37", "16"),
    @("G38", "16"),
    @("D39", "0.02818"),
    @("E39", "11.01%"),
    @("G39", "16"),
    @("D40", "0.05550"),
    @("E40", "5.46%"),
    @("G40", "16"),
    @("D41", "0.007624"),
    @("E41", "4.59%"),
    @("G41", "16"),
    @("D42", "0.009839"),
    @("E42", "9.96%"),
    @("G42", "16"),
    @("D43", "0.1420"),
    @("E43", "-0.92%"),
    @("G43", "16"),
    @("D44", "0.002103"),
    @("E44", "-3.20%"),
    @("G44", "16"),
    @("D45", "0.01182"),
    @("E45", "22.42%"),
    @("G45", "16"),
    @("D46", "0.00006767"),
    @("E46", "0.01%"),
    @("G46", "16"),
    @("D47", "0.00000000749"),
    @("E47", "0.03%"),
    @("G47", "16"),
    @("B48", "CoinbaseStockToken"),
    @("C48", "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"),
    @("D48", "0.002275"),
    @("E48", "26.65%"),
    @("G48", "16"),
    @("B49", "BOLO"),
    @("C49", "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"),
    @("D49", "0.003066"),
    @("E49", "6.93%"),
    @("G49", "16"),
    @("D50", "0.00002096"),
    @("E50", "0.03%"),
    @("G50", "16"),
    @("D51", "0.0001996"),
    @("E51", "0.03%"),
    @("G51", "16")
)

foreach ($update in $cellUpdates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
}

